$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number that Excel would
# otherwise auto-convert from text to a Number (losing formatting such
# as trailing zeros, e.g. "1.000" -> 1). Force text format first so the
# literal string is preserved, matching the source data (openpyxl authored,
# these are inline-string price cells, not real numbers).
$textCells = @(
    @{ Addr = 'D4'; Value = '1.000' }
    @{ Addr = 'D5'; Value = '228.78' }
    @{ Addr = 'D6'; Value = '1.000' }
    @{ Addr = 'D7'; Value = '0.5229' }
    @{ Addr = 'D8'; Value = '0.2743' }
    @{ Addr = 'D9'; Value = '39.27' }
    @{ Addr = 'D12'; Value = '0.07070' }
    @{ Addr = 'D13'; Value = '14.94' }
    @{ Addr = 'D14'; Value = '0.6357' }
    @{ Addr = 'D15'; Value = '4.517' }
    @{ Addr = 'D16'; Value = '76.58' }
    @{ Addr = 'D17'; Value = '1.000' }
    @{ Addr = 'D18'; Value = '1.0000' }
    @{ Addr = 'D20'; Value = '11.47' }
    @{ Addr = 'D21'; Value = '0.000006632' }
    @{ Addr = 'D23'; Value = '4.218' }
    @{ Addr = 'D24'; Value = '8.781' }
    @{ Addr = 'D25'; Value = '5.146' }
    @{ Addr = 'D26'; Value = '140.29' }
    @{ Addr = 'D27'; Value = '1.506' }
    @{ Addr = 'D28'; Value = '15.04' }
    @{ Addr = 'D29'; Value = '1.778' }
    @{ Addr = 'D30'; Value = '102.05' }
    @{ Addr = 'D31'; Value = '0.08296' }
    @{ Addr = 'D32'; Value = '3.708' }
    @{ Addr = 'D33'; Value = '3.506' }
    @{ Addr = 'D34'; Value = '0.04448' }
    @{ Addr = 'D35'; Value = '2.612' }
    @{ Addr = 'D36'; Value = '0.9683' }
    @{ Addr = 'D37'; Value = '0.6172' }
    @{ Addr = 'D38'; Value = '2.668' }
    @{ Addr = 'D39'; Value = '0.01567' }
    @{ Addr = 'D40'; Value = '0.9998' }
    @{ Addr = 'D41'; Value = '1.893' }
    @{ Addr = 'D43'; Value = '0.3810' }
    @{ Addr = 'D44'; Value = '5.017' }
    @{ Addr = 'D45'; Value = '0.7224' }
    @{ Addr = 'D46'; Value = '0.05328' }
    @{ Addr = 'D47'; Value = '0.1115' }
    @{ Addr = 'D50'; Value = '29.92' }
    @{ Addr = 'D51'; Value = '7.538' }
)
foreach ($item in $textCells) {
    $cell = $ws.Range($item.Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Value
}

# Remaining D/E cell updates (already non-numeric-looking text, e.g.
# containing multiple dots, a "%" sign, or surrounding spaces) can be
# assigned directly.
$ws.Range('D2').Value = '25.820.53'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = '1.731.49'
$ws.Range('E3').Value = '  -2.00%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('E5').Value = '  -3.79%  '
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('E8').Value = '  -0.43%  '
$ws.Range('E9').Value = '  -3.25%  '
$ws.Range('E10').Value = '  -1.50%  '
$ws.Range('D11').Value = '1.739.22'
$ws.Range('E11').Value = '  -1.55%  '
$ws.Range('E12').Value = '  +0.27%  '
$ws.Range('E13').Value = '  -6.90%  '
$ws.Range('E14').Value = '  -2.66%  '
$ws.Range('E15').Value = '  -0.21%  '
$ws.Range('E16').Value = '  -2.51%  '
$ws.Range('E17').Value = '  +0.08%  '
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('D19').Value = '25.818.59'
$ws.Range('E19').Value = '  -1.06%  '
$ws.Range('E20').Value = '  -2.28%  '
$ws.Range('E21').Value = '  -1.62%  '
$ws.Range('D22').Value = '1.960.21'
$ws.Range('E22').Value = '  -1.90%  '
$ws.Range('E23').Value = '  +2.88%  '
$ws.Range('E24').Value = '  +4.17%  '
$ws.Range('E25').Value = '  -1.34%  '
$ws.Range('E26').Value = '  +1.66%  '
$ws.Range('E27').Value = '  +1.55%  '
$ws.Range('E28').Value = '  -1.43%  '
$ws.Range('E29').Value = '  -3.81%  '
$ws.Range('E30').Value = '  -1.22%  '
$ws.Range('E31').Value = '  -1.67%  '
$ws.Range('E32').Value = '  -0.84%  '
$ws.Range('E33').Value = '  +1.13%  '
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('E35').Value = '  -1.50%  '
$ws.Range('E36').Value = '  -3.98%  '
$ws.Range('E37').Value = '  +0.52%  '
$ws.Range('E38').Value = '  -3.27%  '
$ws.Range('E39').Value = '  -1.30%  '
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('E41').Value = '  -5.54%  '
$ws.Range('E42').Value = '  -3.52%  '
$ws.Range('E43').Value = '  -2.73%  '
$ws.Range('E44').Value = '  +0.54%  '
$ws.Range('E45').Value = '  -4.21%  '
$ws.Range('E46').Value = '  -3.27%  '
$ws.Range('E47').Value = '  -0.85%  '
$ws.Range('E48').Value = '  -7.46%  '
$ws.Range('E49').Value = '  +0.32%  '
$ws.Range('E50').Value = '  -1.27%  '
$ws.Range('E51').Value = '  +0.92%  '
